$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 28.743868
$ws.Range("H2").Value = 86.231604
$ws.Range("I2").Value = 0.0554303735704667
$ws.Range("J2").Value = 0.0554303735704667
$ws.Range("M2").Value = 0.05871566666666667
$ws.Range("N2").Value = 0.176147
$ws.Range("O2").Value = 0.008355169877275808
$ws.Range("P2").Value = 0.008355169877275808
$ws.Range("Q2").Value = 1.687715372198667
$ws.Range("R2").Value = 15.189438349788
$ws.Range("S2").Value = 0.0004631301875421085
$ws.Range("T2").Value = 0.0004631301875421085
$ws.Range("G3").Value = 28.743868
$ws.Range("H3").Value = 86.231604
$ws.Range("I3").Value = 0.0554303735704667
$ws.Range("J3").Value = 0.0554303735704667
$ws.Range("O3").Value = 0.1868088427899751
$ws.Range("P3").Value = 0.1868088427899751
$ws.Range("Q3").Value = 37.73473912203468
$ws.Range("R3").Value = 339.612652098312
$ws.Range("S3").Value = 0.01035488394211491
$ws.Range("T3").Value = 0.01035488394211491
$ws.Range("G4").Value = 28.743868
$ws.Range("H4").Value = 86.231604
$ws.Range("I4").Value = 0.0554303735704667
$ws.Range("J4").Value = 0.0554303735704667
$ws.Range("O4").Value = 0.8048359873327491
$ws.Range("P4").Value = 0.8048359873327491
$ws.Range("Q4").Value = 162.574081421676
$ws.Range("R4").Value = 1463.166732795084
$ws.Range("S4").Value = 0.04461235944080968
$ws.Range("T4").Value = 0.04461235944080968
$ws.Range("I5").Value = 0.848161237947095
$ws.Range("J5").Value = 0.8481612379470951
$ws.Range("M5").Value = 0.05871566666666667
$ws.Range("N5").Value = 0.176147
$ws.Range("O5").Value = 0.008355169877275808
$ws.Range("P5").Value = 0.008355169877275808
$ws.Range("Q5").Value = 25.82437510666611
$ws.Range("R5").Value = 232.419375959995
$ws.Range("S5").Value = 0.007086531226368527
$ws.Range("T5").Value = 0.007086531226368528
$ws.Range("I6").Value = 0.848161237947095
$ws.Range("J6").Value = 0.8481612379470951
$ws.Range("O6").Value = 0.1868088427899751
$ws.Range("P6").Value = 0.1868088427899751
$ws.Range("S6").Value = 0.1584440193602095
$ws.Range("T6").Value = 0.1584440193602095
$ws.Range("I7").Value = 0.848161237947095
$ws.Range("J7").Value = 0.8481612379470951
$ws.Range("O7").Value = 0.8048359873327491
$ws.Range("P7").Value = 0.8048359873327491
$ws.Range("S7").Value = 0.6826306873605169
$ws.Range("T7").Value = 0.682630687360517
$ws.Range("G8").Value = 49.99334866666666
$ws.Range("I8").Value = 0.09640838848243828
$ws.Range("J8").Value = 0.09640838848243828
$ws.Range("M8").Value = 0.05871566666666667
$ws.Range("N8").Value = 0.176147
$ws.Range("O8").Value = 0.008355169877275808
$ws.Range("P8").Value = 0.008355169877275808
$ws.Range("Q8").Value = 2.935392795862444
$ws.Range("R8").Value = 26.418535162762
$ws.Range("S8").Value = 0.0008055084633651723
$ws.Range("T8").Value = 0.0008055084633651723
$ws.Range("G9").Value = 49.99334866666666
$ws.Range("I9").Value = 0.09640838848243828
$ws.Range("J9").Value = 0.09640838848243828
$ws.Range("O9").Value = 0.1868088427899751
$ws.Range("P9").Value = 0.1868088427899751
$ws.Range("R9").Value = 590.6781136053879
$ws.Range("S9").Value = 0.01800993948765066
$ws.Range("T9").Value = 0.01800993948765066
$ws.Range("G10").Value = 49.99334866666666
$ws.Range("I10").Value = 0.09640838848243828
$ws.Range("J10").Value = 0.09640838848243828
$ws.Range("O10").Value = 0.8048359873327491
$ws.Range("P10").Value = 0.8048359873327491
$ws.Range("S10").Value = 0.07759294053142245
$ws.Range("T10").Value = 0.07759294053142245
